$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value2 = [double]1.334738731384277
$ws.Range("B1").Value2 = [double]1.612576246261597
$ws.Range("C1").Value2 = [double]3.103099584579468
$ws.Range("D1").Value2 = [double]1.511473536491394
$ws.Range("E1").Value2 = [double]0.8279739618301392
